$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-set A67's value (tiny floating point refresh as seen in the source diff)
$ws.Range("A67").Value = 44380.76713137153

# Append new row 68 with freshly retrieved data
$ws.Range("A68").Value = 44381.7671101146
$ws.Range("A68").NumberFormat = $ws.Range("A67").NumberFormat

$ws.Range("B68").Value = 78396
$ws.Range("C68").Value = 66085
$ws.Range("D68").Value = 3433
$ws.Range("E68").Value = 2140
$ws.Range("F68").Value = 1518
$ws.Range("G68").Value = 20844
$ws.Range("H68").Value = 1504
$ws.Range("I68").Value = 878
$ws.Range("J68").Value = 190
